$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I16").Value = 'sv'
$ws.Range("J16").Value = 'Statement-opinion'
$ws.Range("I18").Value = 'sv'
$ws.Range("J18").Value = 'Statement-opinion'
$ws.Range("I27").Value = 'ba'
$ws.Range("J27").Value = 'Appreciation'
$ws.Range("I34").Value = 'ba'
$ws.Range("J34").Value = 'Appreciation'
$ws.Range("I35").Value = 'ba'
$ws.Range("J35").Value = 'Appreciation'
$ws.Range("I43").Value = 'b'
$ws.Range("J43").Value = 'Acknowledge (Backchannel)'
$ws.Range("I61").Value = 'sd'
$ws.Range("J61").Value = 'Statement-non-opinion'
$ws.Range("I82").Value = 'sv'
$ws.Range("J82").Value = 'Statement-opinion'
$ws.Range("I93").Value = 'sd'
$ws.Range("J93").Value = 'Statement-non-opinion'
$ws.Range("I94").Value = 'sd'
$ws.Range("J94").Value = 'Statement-non-opinion'
$ws.Range("I101").Value = '%'
$ws.Range("J101").Value = 'Uninterpretable'
$ws.Range("I103").Value = 'sd'
$ws.Range("J103").Value = 'Statement-non-opinion'
$ws.Range("I109").Value = 'sd'
$ws.Range("J109").Value = 'Statement-non-opinion'
$ws.Range("I112").Value = 'sd'
$ws.Range("J112").Value = 'Statement-non-opinion'
$ws.Range("I126").Value = 'aa'
$ws.Range("J126").Value = 'Agree/Accept'
$ws.Range("I127").Value = 'ba'
$ws.Range("J127").Value = 'Appreciation'
$ws.Range("I131").Value = 'sd'
$ws.Range("J131").Value = 'Statement-non-opinion'
$ws.Range("I133").Value = 'aa'
$ws.Range("J133").Value = 'Agree/Accept'
$ws.Range("I139").Value = 'sv'
$ws.Range("J139").Value = 'Statement-opinion'
$ws.Range("I140").Value = 'qy'
$ws.Range("J140").Value = 'Yes-No-Question'
$ws.Range("I141").Value = 'sd'
$ws.Range("J141").Value = 'Statement-non-opinion'
$ws.Range("I152").Value = 'sd'
$ws.Range("J152").Value = 'Statement-non-opinion'
$ws.Range("I154").Value = 'sd'
$ws.Range("J154").Value = 'Statement-non-opinion'
$ws.Range("I155").Value = 'sd'
$ws.Range("J155").Value = 'Statement-non-opinion'
$ws.Range("I158").Value = 'sd'
$ws.Range("J158").Value = 'Statement-non-opinion'
$ws.Range("I162").Value = 'sd'
$ws.Range("J162").Value = 'Statement-non-opinion'
$ws.Range("I163").Value = 'sd'
$ws.Range("J163").Value = 'Statement-non-opinion'
$ws.Range("I168").Value = 'sv'
$ws.Range("J168").Value = 'Statement-opinion'
$ws.Range("I207").Value = 'sd'
$ws.Range("J207").Value = 'Statement-non-opinion'
$ws.Range("I215").Value = 'sv'
$ws.Range("J215").Value = 'Statement-opinion'
$ws.Range("I218").Value = 'sv'
$ws.Range("J218").Value = 'Statement-opinion'
$ws.Range("I245").Value = 'ba'
$ws.Range("J245").Value = 'Appreciation'
$ws.Range("I246").Value = 'b'
$ws.Range("J246").Value = 'Acknowledge (Backchannel)'
$ws.Range("I250").Value = 'sd'
$ws.Range("J250").Value = 'Statement-non-opinion'
$ws.Range("I264").Value = 'sd'
$ws.Range("J264").Value = 'Statement-non-opinion'
$ws.Range("I272").Value = 'sd'
$ws.Range("J272").Value = 'Statement-non-opinion'
$ws.Range("I273").Value = 'sd'
$ws.Range("J273").Value = 'Statement-non-opinion'
$ws.Range("I280").Value = 'sv'
$ws.Range("J280").Value = 'Statement-opinion'
$ws.Range("I286").Value = 'aa'
$ws.Range("J286").Value = 'Agree/Accept'
$ws.Range("I293").Value = 'sd'
$ws.Range("J293").Value = 'Statement-non-opinion'
$ws.Range("I297").Value = '%'
$ws.Range("J297").Value = 'Uninterpretable'
$ws.Range("I314").Value = 'sv'
$ws.Range("J314").Value = 'Statement-opinion'
$ws.Range("I316").Value = 'sd'
$ws.Range("J316").Value = 'Statement-non-opinion'
$ws.Range("I323").Value = 'b'
$ws.Range("J323").Value = 'Acknowledge (Backchannel)'
$ws.Range("I326").Value = 'sv'
$ws.Range("J326").Value = 'Statement-opinion'
$ws.Range("I328").Value = 'sv'
$ws.Range("J328").Value = 'Statement-opinion'
$ws.Range("I339").Value = 'sd'
$ws.Range("J339").Value = 'Statement-non-opinion'
$ws.Range("I345").Value = 'sv'
$ws.Range("J345").Value = 'Statement-opinion'
$ws.Range("I351").Value = 'sd'
$ws.Range("J351").Value = 'Statement-non-opinion'
$ws.Range("I352").Value = 'sd'
$ws.Range("J352").Value = 'Statement-non-opinion'
$ws.Range("I353").Value = '%'
$ws.Range("J353").Value = 'Uninterpretable'
$ws.Range("I358").Value = 'ba'
$ws.Range("J358").Value = 'Appreciation'
$ws.Range("I359").Value = 'sd'
$ws.Range("J359").Value = 'Statement-non-opinion'
$ws.Range("I370").Value = 'sd'
$ws.Range("J370").Value = 'Statement-non-opinion'
$ws.Range("I387").Value = 'sd'
$ws.Range("J387").Value = 'Statement-non-opinion'
$ws.Range("I390").Value = 'b'
$ws.Range("J390").Value = 'Acknowledge (Backchannel)'
$ws.Range("I392").Value = 'sd'
$ws.Range("J392").Value = 'Statement-non-opinion'
$ws.Range("I406").Value = 'aa'
$ws.Range("J406").Value = 'Agree/Accept'
$ws.Range("I412").Value = 'aa'
$ws.Range("J412").Value = 'Agree/Accept'
$ws.Range("I444").Value = 'sv'
$ws.Range("J444").Value = 'Statement-opinion'
$ws.Range("I453").Value = 'sd'
$ws.Range("J453").Value = 'Statement-non-opinion'
$ws.Range("I458").Value = '%'
$ws.Range("J458").Value = 'Uninterpretable'
$ws.Range("I469").Value = 'sd'
$ws.Range("J469").Value = 'Statement-non-opinion'
$ws.Range("I477").Value = 'sd'
$ws.Range("J477").Value = 'Statement-non-opinion'
$ws.Range("I478").Value = 'ba'
$ws.Range("J478").Value = 'Appreciation'
$ws.Range("I499").Value = 'sd'
$ws.Range("J499").Value = 'Statement-non-opinion'
$ws.Range("I501").Value = 'b'
$ws.Range("J501").Value = 'Acknowledge (Backchannel)'
$ws.Range("I521").Value = 'aa'
$ws.Range("J521").Value = 'Agree/Accept'
$ws.Range("I522").Value = 'aa'
$ws.Range("J522").Value = 'Agree/Accept'
$ws.Range("I527").Value = 'sv'
$ws.Range("J527").Value = 'Statement-opinion'
$ws.Range("I535").Value = 'sd'
$ws.Range("J535").Value = 'Statement-non-opinion'
$ws.Range("I536").Value = 'sd'
$ws.Range("J536").Value = 'Statement-non-opinion'
$ws.Range("I537").Value = 'sd'
$ws.Range("J537").Value = 'Statement-non-opinion'
$ws.Range("I542").Value = 'ba'
$ws.Range("J542").Value = 'Appreciation'
$ws.Range("I551").Value = 'b'
$ws.Range("J551").Value = 'Acknowledge (Backchannel)'
$ws.Range("I552").Value = 'sv'
$ws.Range("J552").Value = 'Statement-opinion'
$ws.Range("I555").Value = 'sd'
$ws.Range("J555").Value = 'Statement-non-opinion'
$ws.Range("I556").Value = 'b'
$ws.Range("J556").Value = 'Acknowledge (Backchannel)'
$ws.Range("I557").Value = 'aa'
$ws.Range("J557").Value = 'Agree/Accept'
$ws.Range("I561").Value = 'sd'
$ws.Range("J561").Value = 'Statement-non-opinion'
$ws.Range("I566").Value = 'sv'
$ws.Range("J566").Value = 'Statement-opinion'
$ws.Range("I589").Value = 'sd'
$ws.Range("J589").Value = 'Statement-non-opinion'
$ws.Range("I592").Value = 'sd'
$ws.Range("J592").Value = 'Statement-non-opinion'
$ws.Range("I604").Value = 'aa'
$ws.Range("J604").Value = 'Agree/Accept'
$ws.Range("I607").Value = 'sd'
$ws.Range("J607").Value = 'Statement-non-opinion'
$ws.Range("I609").Value = 'sv'
$ws.Range("J609").Value = 'Statement-opinion'
$ws.Range("I611").Value = 'sd'
$ws.Range("J611").Value = 'Statement-non-opinion'
$ws.Range("I615").Value = 'ba'
$ws.Range("J615").Value = 'Appreciation'
$ws.Range("I619").Value = 'aa'
$ws.Range("J619").Value = 'Agree/Accept'
$ws.Range("I630").Value = 'aa'
$ws.Range("J630").Value = 'Agree/Accept'
$ws.Range("I644").Value = 'b'
$ws.Range("J644").Value = 'Acknowledge (Backchannel)'
$ws.Range("I665").Value = 'b'
$ws.Range("J665").Value = 'Acknowledge (Backchannel)'
$ws.Range("I684").Value = 'sd'
$ws.Range("J684").Value = 'Statement-non-opinion'
$ws.Range("I685").Value = 'aa'
$ws.Range("J685").Value = 'Agree/Accept'
$ws.Range("I687").Value = 'ba'
$ws.Range("J687").Value = 'Appreciation'
